# Add data for 2022-03-18: advance the "through" date from 03-09 to 03-10,
# and update the totals that changed as a result (March + grand total).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-03-10"

# Update the 2022 column header text (shared string used by I1).
$ws.Range("I1").Value = "2022 (through 03-10)"

# Update the March 2022 figure (row 4 = March).
$ws.Range("I4").Value = 45

# Update the Total 2022 figure (row 14 = Total).
$ws.Range("I14").Value = 346
